# React 101 Homework updated.
# Fill in the newly-added quiz/homework columns (C:H) with scores for
# Q02, Q03, H01, H02, H03, H04 and update the Project (C column) scores
# in the summary table further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 headers: new assessment columns ---
$ws.Range("C1").Value = "Q02"
$ws.Range("D1").Value = "Q03"
$ws.Range("E1").Value = "H01"
$ws.Range("F1").Value = "H02"
$ws.Range("G1").Value = "H03"
$ws.Range("H1").Value = "H04"

# --- Row 2 (student 1) ---
$ws.Range("C2").Formula = "=(10/15)*10"
$ws.Range("D2").Formula = "=(15/25)*10"
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0

# --- Row 3 (student 2) ---
$ws.Range("C3").Formula = "=(11/15)*10"
$ws.Range("D3").Formula = "=(23/25)*10"
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 9.5
$ws.Range("G3").Value = 9
$ws.Range("H3").Value = 9.75

# --- Row 4 (student 3) ---
$ws.Range("C4").Formula = "=(5/15)*10"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = 9.5
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0

# --- Row 5 (student 4) ---
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 7
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0

# --- Row 6 (student 5) ---
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0

# --- Row 7 (student 6) ---
$ws.Range("C7").Formula = "=(8/15)*10"
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 9
$ws.Range("F7").Value = 9.5
$ws.Range("G7").Value = 9
$ws.Range("H7").Value = 8.5

# --- Row 8 (student 7) ---
$ws.Range("C8").Formula = "=(13/15)*10"
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 9.5
$ws.Range("F8").Value = 9.5
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0

# --- Row 9 (student 8) ---
$ws.Range("C9").Formula = "=(6/15)*10"
$ws.Range("D9").Formula = "=(18/25)*10"
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0

# --- Row 10 (student 9) ---
$ws.Range("C10").Formula = "=(10/15)*10"
$ws.Range("D10").Formula = "=(21/25)*10"
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 9.5
$ws.Range("G10").Value = 10
$ws.Range("H10").Value = 8.5

# --- Row 11 (student 10) ---
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0

# --- Row 12 (student 11) ---
$ws.Range("C12").Formula = "=(13/15)*10"
$ws.Range("D12").Formula = "=(21/25)*10"
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 9
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0

# --- Summary table: Project (C) scores, rows 15-25 ---
$ws.Range("C15").Value = 8
$ws.Range("C16").Value = 9.5
$ws.Range("C17").Value = 8.5
$ws.Range("C18").Value = 6
$ws.Range("C19").Value = 6
$ws.Range("C20").Value = 9
$ws.Range("C21").Value = 10
$ws.Range("C22").Value = 7.5
$ws.Range("C23").Value = 8
$ws.Range("C24").Value = 6
$ws.Range("C25").Value = 9

# --- Selection moved by the author while editing ---
$ws.Range("C26").Select()
